# Applies the commit "Testing PivotalTracker [finished #90278336]":
#   - Collapses runs that were split around spell-check (proofErr) markers
#     back into single contiguous runs. Word's Find/Replace operates on
#     rendered text regardless of run boundaries, so re-"finding and
#     replacing" each sentence with its own literal text rewrites it as a
#     single run and drops the now-pointless <w:proofErr/> spell-check
#     bookmarks around "myhours", "vs", "GitHub", "Kaban", "NUnit" and
#     "algos".
#   - Adds a new, empty bullet-list paragraph (ListParagraph style, same
#     numbered list as the "Questions" section) after the final
#     "Share repo now..." paragraph.

$d = $word.ActiveDocument

function Merge-ParagraphRuns($text) {
    $r = $d.Content
    $r.Find.ClearFormatting()
    $r.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, $text, 2) | Out-Null
}

Merge-ParagraphRuns "Treat as client.  New client in myhours and track time.  "
Merge-ParagraphRuns "TFS vs GIT - TFS - choose Scrum process template vs David saying that moving over to GIT - so create a new GitHub repo.  "
Merge-ParagraphRuns "Need to look at Kaban process if using GIT."
Merge-ParagraphRuns "Use NUnit testing - create new test project.  Which runner?"
Merge-ParagraphRuns "Impact on memory key requirement - use ANTs Memory profiler?  Most significant impact will be on the data structure and search algos."

# Add a new empty list paragraph right after "Share repo now so that we can
# discuss if required?" (the last paragraph in the document). The new
# paragraph inherits the ListParagraph style / numId=2 bullet automatically
# from the paragraph it follows.
$anchor = $d.Content
$anchor.Find.ClearFormatting()
$anchor.Find.Execute("Share repo now so that we can discuss if required?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$anchorPara = $anchor.Paragraphs(1)
$anchorPara.Range.InsertParagraphAfter() | Out-Null
